$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Insert a brand-new "2022-Q1" sheet right before the "总计" (totals) sheet.
# ---------------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")
$q1Sheet = $wb.Worksheets.Add($totalSheet)
$q1Sheet.Name = "2022-Q1"

# NOTE: `Worksheets.Add` shifts sheet positions, which can invalidate any
# sheet handle obtained *before* the call (handles here are position based).
# Re-resolve "总计" by name now that the new sheet has been inserted so the
# rest of the script references the correct worksheet.
$totalSheet = $wb.Worksheets.Item("总计")

# Pull header-row (B1:H1) formatting (bold/centered/bordered style) from the
# "2021-Q3" sheet, which already uses the desired look.
$srcSheet = $wb.Worksheets.Item("2021-Q3")
$srcSheet.Range("B1:H1").Copy()
$q1Sheet.Range("B1:H1").PasteSpecial(-4122)

# Pull the A-column index-cell style too (used on A2/A3 in the new sheet).
$srcSheet.Range("A2").Copy()
$q1Sheet.Range("A2:A3").PasteSpecial(-4122)

# Header labels
$q1Sheet.Range("B1").Value = "基金代码"
$q1Sheet.Range("C1").Value = "基金名称"
$q1Sheet.Range("D1").Value = "基金规模"
$q1Sheet.Range("E1").Value = "股票总仓位"
$q1Sheet.Range("F1").Value = "仓位占比"
$q1Sheet.Range("G1").Value = "持有市值(亿元)"
$q1Sheet.Range("H1").Value = "仓位排名"

# Force the B:G data cells to be stored as text (keeps leading zeros on fund
# codes like "090019" and preserves the numeric strings verbatim).
$q1Sheet.Range("B2:G3").NumberFormat = "@"

# Row 2 - 090019 / 大成景恒混合A
$q1Sheet.Range("A2").Value = 0
$q1Sheet.Range("B2").Value = "090019"
$q1Sheet.Range("C2").Value = "大成景恒混合A"
$q1Sheet.Range("D2").Value = "2.31"
$q1Sheet.Range("E2").Value = "93.51"
$q1Sheet.Range("F2").Value = "2.02"
$q1Sheet.Range("G2").Value = "0.0467"
$q1Sheet.Range("H2").Value = 3

# Row 3 - 006038 / 大成景恒混合C
$q1Sheet.Range("A3").Value = 1
$q1Sheet.Range("B3").Value = "006038"
$q1Sheet.Range("C3").Value = "大成景恒混合C"
$q1Sheet.Range("D3").Value = "0.92"
$q1Sheet.Range("E3").Value = "93.51"
$q1Sheet.Range("F3").Value = "2.02"
$q1Sheet.Range("G3").Value = "0.0186"
$q1Sheet.Range("H3").Value = 3

# ---------------------------------------------------------------------------
# 2) Add a new top data row to "总计" for 2022-Q1, pushing the older rows
#    down and renumbering the running index in column A.
# ---------------------------------------------------------------------------
$totalSheet.Rows.Item(2).Insert()
$totalSheet.Range("B2:D2").ClearFormats()

# Give A2 the same index-column style (s=2) used by the other A-cells.
$totalSheet.Range("A3").Copy()
$totalSheet.Range("A2").PasteSpecial(-4122)

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q1"
$totalSheet.Range("C2").Value = 2
$totalSheet.Range("D2").Value = 0.07000000000000001

# Renumber the index column for the rows that shifted down.
$totalSheet.Range("A3").Value = 1
$totalSheet.Range("A4").Value = 2
$totalSheet.Range("A5").Value = 3
$totalSheet.Range("A6").Value = 4

Write-Output "2022-Q1 sheet added and 总计 updated"
